# ------------------------------------------------------------------
# Applies the "output generated at 456a3b4" update to 上海-漫展信息.xlsx
#   * bumps a bunch of "want to go" (F column) counters across all 4
#     sheets
#   * on sheet "演出" (Performances) inserts one new row (a newly
#     scraped concert, 茅原实里动漫交响音乐会) above the existing
#     2024-04-20 "Laurent Coulondre" row, pushing that row and the
#     four rows below it down by one
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ============================================================
# Sheet 1: 展览 (Exhibitions) -- simple F-column counter bumps
# ============================================================
$ws = $wb.Worksheets.Item(1)
$ws.Range("F4").Value = 120
$ws.Range("F5").Value = 343
$ws.Range("F6").Value = 739
$ws.Range("F7").Value = 190
$ws.Range("F8").Value = 233
$ws.Range("F9").Value = 628
$ws.Range("F12").Value = 589
$ws.Range("F13").Value = 491
$ws.Range("F14").Value = 134
$ws.Range("F16").Value = 150
$ws.Range("F17").Value = 804
$ws.Range("F23").Value = 196
$ws.Range("F25").Value = 143
$ws.Range("F26").Value = 579
$ws.Range("F28").Value = 53
$ws.Range("F29").Value = 194
$ws.Range("F30").Value = 584
$ws.Range("F34").Value = 267

# ============================================================
# Sheet 2: 演出 (Performances)
# ============================================================
$ws = $wb.Worksheets.Item(2)

# ---- F-column counter bumps on rows that aren't moving ----
$ws.Range("F4").Value = 1002
$ws.Range("F5").Value = 1002
$ws.Range("F14").Value = 550
$ws.Range("F15").Value = 84
$ws.Range("F19").Value = 32
$ws.Range("F24").Value = 273
$ws.Range("F25").Value = 240

# ---- Insert a new row above row 30 ----
# The engine's Range.Value setter runs new literals through Excel's
# usual "does this look like a date/number" inference, which would
# turn a freshly-typed "2024-04-20" into a date serial. To keep every
# moved cell byte-identical to a real row-insert (plain shared-string
# text, no incidental number formatting), shift the B:I payload of
# rows 30-34 down into rows 31-35 with Copy/PasteSpecial (which
# carries the existing text cells verbatim) instead of re-typing them,
# working from the bottom row up so each source row is read before it
# gets overwritten. Column A already holds the simple 0-based
# row-index counter and is left untouched except for the new row 35.
for ($r = 34; $r -ge 30; $r--) {
    $nr = $r + 1
    $ws.Range("B" + $r + ":I" + $r).Copy()
    $ws.Range("B" + $nr + ":I" + $nr).PasteSpecial()
}
$excel.CutCopyMode = $false

# Extend the A-column index counter onto the newly created row 35
# (value 34, matching the existing s="1" bold/bordered index style).
$ws.Cells.Item(35, 1).Value = 34
$ws.Range("A34").Copy()
$ws.Range("A35").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Fill row 30 with the new concert's data.
$ws.Range("C30").Value = "上海· 茅原实里动漫交响音乐会"
$ws.Range("D30").Value = "东大名路889号 友邦大剧院"
$ws.Range("E30").Value = "2024.04.20 19:30-04.20 21:00"
$ws.Range("F30").Value = 1
$ws.Range("G30").Value = 280
$ws.Range("H30").Value = "https://show.bilibili.com/platform/detail.html?id=81703"
$ws.Range("I30").Value = "//i2.hdslb.com/bfs/openplatform/202402/yiVaqJVK1707016321221.jpeg"
# B30's target text ("2024-04-20") is identical to the now-shifted
# B31, so clone that cell instead of typing the literal again -- same
# reasoning as above, avoids the date-inference side effect.
$ws.Range("B31").Copy()
$ws.Range("B30").PasteSpecial()
$excel.CutCopyMode = $false

# ============================================================
# Sheet 3: 本地生活 (Local life) -- simple F-column counter bumps
# ============================================================
$ws = $wb.Worksheets.Item(3)
$ws.Range("F5").Value = 2374
$ws.Range("F9").Value = 1212
$ws.Range("F10").Value = 317

# ============================================================
# Sheet 4: 全部类型 (All types) -- simple F-column counter bumps
# ============================================================
$ws = $wb.Worksheets.Item(4)
$ws.Range("F5").Value = 2374
$ws.Range("F9").Value = 1212
$ws.Range("F10").Value = 317
$ws.Range("F12").Value = 120
$ws.Range("F13").Value = 343
$ws.Range("F14").Value = 739
$ws.Range("F15").Value = 190
$ws.Range("F17").Value = 233
$ws.Range("F18").Value = 632
$ws.Range("F20").Value = 589
$ws.Range("F21").Value = 1002
$ws.Range("F22").Value = 491
$ws.Range("F23").Value = 134
$ws.Range("F25").Value = 150
$ws.Range("F26").Value = 804
$ws.Range("F31").Value = 196
$ws.Range("F32").Value = 143
$ws.Range("F33").Value = 579
$ws.Range("F35").Value = 550
$ws.Range("F36").Value = 84
$ws.Range("F37").Value = 53
$ws.Range("F38").Value = 194
$ws.Range("F40").Value = 32
$ws.Range("F43").Value = 273
$ws.Range("F44").Value = 273
$ws.Range("F45").Value = 240
$ws.Range("F46").Value = 588
$ws.Range("F50").Value = 267
